# Add two new booking rows (4 and 5) to Sheet1, mirroring the existing
# mock data layout (columns A-K), as part of the updated mock data for
# the React event booking / Excel import feature.
#
# Note: D (Phone) and H (Date) values look numeric/date-like, so a
# leading apostrophe is used to force them to stay as text, matching
# the rest of the sheet's "numbers stored as text" convention.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Test User 11
$ws.Range("A4").Value = 246
$ws.Range("B4").Value = "Test User 11"
$ws.Range("D4").Value = "'1234567890"
$ws.Range("H4").Value = "'2025-12-04"
$ws.Range("I4").Value = "10:45 am - 11:45 am"
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 1

# Row 5: Test User 12
$ws.Range("A5").Value = 247
$ws.Range("B5").Value = "Test User 12"
$ws.Range("D5").Value = "'0987654321"
$ws.Range("H5").Value = "'2025-12-04"
$ws.Range("I5").Value = "10:45 am - 11:45 am"
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 1

Write-Output "Added rows 4 and 5 with new test user booking data."
